# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 19:46"

# Swap the display order of "La Palma" / "Arroyo de la Luz" rows (A56 <-> A57)
$ws.Range("A56").Value = "Arroyo de la Luz"
$ws.Range("A57").Value = "La Palma"

# Update "Casos activos" (column C) for several Galicia provinces
$ws.Range("C18").Value = 5   # A Coruña
$ws.Range("C28").Value = 5   # Pontevedra
$ws.Range("C44").Value = 5   # Ourense
$ws.Range("C48").Value = 5   # Lugo
